$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1498.8
$ws.Range("I28").Value = 1498
$ws.Range("K28").Value = 1498
$ws.Range("M28").Value = -1013
$ws.Range("H49").Value = 2725.2
$ws.Range("I49").Value = 1209
$ws.Range("J49").Value = 4999.5
$ws.Range("K49").Value = 3627
$ws.Range("L49").Value = 14998.5
$ws.Range("M49").Value = -3491
$ws.Range("N49").Value = -15270.5
$ws.Range("H62").Value = 1577.8
$ws.Range("I62").Value = 1548.5
$ws.Range("K62").Value = 1548.5
$ws.Range("M62").Value = -924.5
$ws.Range("H64").Value = 22731408
$ws.Range("I64").Value = 25004048
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 25004048
$ws.Range("L64").Value = 5000
$ws.Range("M64").Value = -25003800
$ws.Range("N64").Value = -5496
$ws.Range("H65").Value = 1577.8
$ws.Range("I65").Value = 1548.5
$ws.Range("K65").Value = 7742.5
$ws.Range("M65").Value = -4622.5
$ws.Range("H67").Value = 22731408
$ws.Range("I67").Value = 25004048
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 25004048
$ws.Range("L67").Value = 5000
$ws.Range("M67").Value = -25003190
$ws.Range("N67").Value = -6716
$ws.Range("H94").Value = 2345.1667
$ws.Range("I94").Value = 2345.1667
$ws.Range("K94").Value = 2345.1667
$ws.Range("M94").Value = -1894.1667
$ws.Range("H111").Value = 1819.75
$ws.Range("I111").Value = 889.5
$ws.Range("J111").Value = 2750
$ws.Range("K111").Value = 2668.5
$ws.Range("L111").Value = 8250
$ws.Range("M111").Value = 398.5
$ws.Range("N111").Value = -14384
$ws.Range("H113").Value = 3800.077
$ws.Range("I113").Value = 4379.625
$ws.Range("J113").Value = 2872.8
$ws.Range("K113").Value = 4379.625
$ws.Range("L113").Value = 2872.8
$ws.Range("M113").Value = -1125.625
$ws.Range("N113").Value = -9380.799999999999
$ws.Range("H132").Value = 1732.7188
$ws.Range("I132").Value = 1764.9333
$ws.Range("K132").Value = 5294.7999
$ws.Range("M132").Value = -2764.7999

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3753.6428
$ws.Range("I32").Value = 2621.8718
$ws.Range("K32").Value = 2621.8718
$ws.Range("M32").Value = -2334.8718
$ws.Range("H45").Value = 6082.8823
$ws.Range("I45").Value = 7490.727
$ws.Range("K45").Value = 7490.727
$ws.Range("M45").Value = -7113.727
$ws.Range("H46").Value = 9995
$ws.Range("J46").Value = 9995
$ws.Range("L46").Value = 9995
$ws.Range("N46").Value = -10633
$ws.Range("H122").Value = 6317.357
$ws.Range("I122").Value = 5743.778
$ws.Range("J122").Value = 7349.8
$ws.Range("K122").Value = 17231.334
$ws.Range("L122").Value = 22049.4
$ws.Range("M122").Value = -14781.334
$ws.Range("N122").Value = -26949.4
$ws.Range("H123").Value = 58984.5
$ws.Range("J123").Value = 58984.5
$ws.Range("L123").Value = 58984.5
$ws.Range("N123").Value = -68784.5
$ws.Range("H132").Value = 3032427.8
$ws.Range("I132").Value = 3449978.2
$ws.Range("J132").Value = 5187
$ws.Range("K132").Value = 10349934.6
$ws.Range("L132").Value = 15561
$ws.Range("M132").Value = -10347404.6
$ws.Range("N132").Value = -20621

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 2517.913
$ws.Range("I94").Value = 2508.9412
$ws.Range("J94").Value = 2543.3333
$ws.Range("K94").Value = 2508.9412
$ws.Range("L94").Value = 2543.3333
$ws.Range("M94").Value = -2057.9412
$ws.Range("N94").Value = -3445.3333

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1567.75
$ws.Range("J19").Value = 1698.3334
$ws.Range("L19").Value = 1698.3334
$ws.Range("N19").Value = -2038.3334
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H24").Value = 1567.75
$ws.Range("J24").Value = 1698.3334
$ws.Range("L24").Value = 1698.3334
$ws.Range("N24").Value = -2038.3334
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H41").Value = 20160
$ws.Range("I41").Value = 5400
$ws.Range("J41").Value = 30000
$ws.Range("K41").Value = 5400
$ws.Range("L41").Value = 30000
$ws.Range("M41").Value = -4972
$ws.Range("N41").Value = -30856
$ws.Range("H88").Value = 17486.75
$ws.Range("J88").Value = 19973.5
$ws.Range("L88").Value = 19973.5
$ws.Range("N88").Value = -20785.5
$ws.Range("H91").Value = 17486.75
$ws.Range("J91").Value = 19973.5
$ws.Range("L91").Value = 19973.5
$ws.Range("N91").Value = -22781.5
$ws.Range("H122").Value = 2709
$ws.Range("I122").Value = 2709
$ws.Range("K122").Value = 8127
$ws.Range("M122").Value = -5677
$ws.Range("H129").Value = 86000
$ws.Range("J129").Value = 89000
$ws.Range("L129").Value = 89000
$ws.Range("N129").Value = -99000

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 135185.95
$ws.Range("I11").Value = 137037.84
$ws.Range("K11").Value = 411113.52
$ws.Range("M11").Value = -410973.52
$ws.Range("H23").Value = 40.22222
$ws.Range("J23").Value = 41.2
$ws.Range("L23").Value = 123.6
$ws.Range("N23").Value = -593.6
$ws.Range("H32").Value = 75275016
$ws.Range("J32").Value = 69
$ws.Range("L32").Value = 207
$ws.Range("N32").Value = -773
$ws.Range("H103").Value = 1158.5385
$ws.Range("I103").Value = 259.75
$ws.Range("J103").Value = 2596.6
$ws.Range("K103").Value = 779.25
$ws.Range("L103").Value = 7789.799999999999
$ws.Range("M103").Value = 99.75
$ws.Range("N103").Value = -9547.799999999999
$ws.Range("H131").Value = 2708.3333
$ws.Range("I131").Value = 2091
$ws.Range("J131").Value = 3943
$ws.Range("K131").Value = 6273
$ws.Range("L131").Value = 11829
$ws.Range("M131").Value = -1233
$ws.Range("N131").Value = -21909
$ws.Range("H132").Value = 4130.875
$ws.Range("I132").Value = 3266.3333
$ws.Range("K132").Value = 29396.9997
$ws.Range("M132").Value = -26866.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 100000
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 100000
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 100000
$ws.Range("N5").Value = -100224
$ws.Range("M5").ClearContents()
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()
$ws.Range("H31").Value = 1714
$ws.Range("I31").Value = 1483.75
$ws.Range("K31").Value = 1483.75
$ws.Range("M31").Value = -1191.75
$ws.Range("H37").Value = 1714
$ws.Range("I37").Value = 1483.75
$ws.Range("K37").Value = 1483.75
$ws.Range("M37").Value = -1206.75
$ws.Range("H50").Value = 0
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()
$ws.Range("H70").Value = 6262.625
$ws.Range("I70").Value = 6125.75
$ws.Range("J70").Value = 6399.5
$ws.Range("K70").Value = 6125.75
$ws.Range("L70").Value = 6399.5
$ws.Range("M70").Value = -5855.75
$ws.Range("N70").Value = -6939.5
$ws.Range("H73").Value = 6262.625
$ws.Range("I73").Value = 6125.75
$ws.Range("J73").Value = 6399.5
$ws.Range("K73").Value = 6125.75
$ws.Range("L73").Value = 6399.5
$ws.Range("M73").Value = -5189.75
$ws.Range("N73").Value = -8271.5
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1790.375
$ws.Range("I46").Value = 1790.375
$ws.Range("K46").Value = 1790.375
$ws.Range("M46").Value = -1602.375
$ws.Range("H82").Value = 2105.7368
$ws.Range("I82").Value = 2059.6
$ws.Range("J82").Value = 2157
$ws.Range("K82").Value = 2059.6
$ws.Range("L82").Value = 2157
$ws.Range("M82").Value = -1698.6
$ws.Range("N82").Value = -2879
$ws.Range("H85").Value = 2105.7368
$ws.Range("I85").Value = 2059.6
$ws.Range("J85").Value = 2157
$ws.Range("K85").Value = 2059.6
$ws.Range("L85").Value = 2157
$ws.Range("M85").Value = -811.5999999999999
$ws.Range("N85").Value = -4653
$ws.Range("H132").Value = 36926476
$ws.Range("J132").Value = 2999
$ws.Range("L132").Value = 8997
$ws.Range("N132").Value = -14057
$ws.Range("H136").Value = 1834.9333
$ws.Range("I136").Value = 1834.9333
$ws.Range("K136").Value = 5504.7999
$ws.Range("M136").Value = -2954.7999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4999.6665
$ws.Range("I62").Value = 4999.6665
$ws.Range("K62").Value = 4999.6665
$ws.Range("M62").Value = -4375.6665
$ws.Range("H65").Value = 4999.6665
$ws.Range("I65").Value = 4999.6665
$ws.Range("K65").Value = 24998.3325
$ws.Range("M65").Value = -21878.3325
$ws.Range("H107").Value = 507.875
$ws.Range("I107").Value = 406.33334
$ws.Range("J107").Value = 812.5
$ws.Range("K107").Value = 1219.00002
$ws.Range("L107").Value = 2437.5
$ws.Range("M107").Value = 700.9999800000001
$ws.Range("N107").Value = -6277.5
$ws.Range("H113").Value = 854.7143
$ws.Range("I113").Value = 848.03845
$ws.Range("J113").Value = 874
$ws.Range("K113").Value = 2544.11535
$ws.Range("L113").Value = 2622
$ws.Range("M113").Value = -374.11535
$ws.Range("N113").Value = -6962
$ws.Range("H126").Value = 2583.9
$ws.Range("I126").Value = 2676.6667
$ws.Range("J126").Value = 1749
$ws.Range("K126").Value = 8030.000100000001
$ws.Range("L126").Value = 5247
$ws.Range("M126").Value = -5560.000100000001
$ws.Range("N126").Value = -10187
